$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was recorded. Insert a row right after the
# header's first data row (row 2) - i.e. at row 3 - shifting every
# existing record down by one, then fill the new row with the new data.
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = 10
$ws.Range("B3").Value = "Vega Modelo de Temuco"
$ws.Range("C3").Value = "La Araucanía"
$ws.Range("D3").Value = 44921
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100103
$ws.Range("H3").Value = "Frutos de hueso (carozo)"
$ws.Range("I3").Value = 100103003
$ws.Range("J3").Value = "Damasco"
$ws.Range("K3").Value = "Dina"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 45
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 20000
$ws.Range("Q3").Value = "$/caja 15 kilos"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1333
$ws.Range("T3").Value = 15
